$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values, preserving original cell style/number format
# by forcing text interpretation only while assigning the literal value.
$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.203.10'
$cell.Style = $origStyle

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.660.24'
$cell.Style = $origStyle

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = $origStyle

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '216.79'
$cell.Style = $origStyle

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5151'
$cell.Style = $origStyle

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.2646'
$cell.Style = $origStyle

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06273'
$cell.Style = $origStyle

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '20.81'
$cell.Style = $origStyle

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.07752'
$cell.Style = $origStyle

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.473'
$cell.Style = $origStyle

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.665.01'
$cell.Style = $origStyle

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.887.68'
$cell.Style = $origStyle

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5452'
$cell.Style = $origStyle

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0₅8106'
$cell.Style = $origStyle

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '64.91'
$cell.Style = $origStyle

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.215.65'
$cell.Style = $origStyle

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.618'
$cell.Style = $origStyle

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '192.62'
$cell.Style = $origStyle

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.09'
$cell.Style = $origStyle

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.009'
$cell.Style = $origStyle

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '139.90'
$cell.Style = $origStyle

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.1223'
$cell.Style = $origStyle

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.269'
$cell.Style = $origStyle

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '16.16'
$cell.Style = $origStyle

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.436'
$cell.Style = $origStyle

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.05974'
$cell.Style = $origStyle

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.272'
$cell.Style = $origStyle

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.572'
$cell.Style = $origStyle

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.259'
$cell.Style = $origStyle

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.594'
$cell.Style = $origStyle

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.9670'
$cell.Style = $origStyle

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.425'
$cell.Style = $origStyle

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.767'
$cell.Style = $origStyle

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5693'
$cell.Style = $origStyle

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.036'
$cell.Style = $origStyle

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.01595'
$cell.Style = $origStyle

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.8566'
$cell.Style = $origStyle

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.011.21'
$cell.Style = $origStyle

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '100.31'
$cell.Style = $origStyle

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.801.16'
$cell.Style = $origStyle

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0₈109'
$cell.Style = $origStyle

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '56.73'
$cell.Style = $origStyle

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.049'
$cell.Style = $origStyle

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.05168'
$cell.Style = $origStyle

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.447'
$cell.Style = $origStyle

# Update Volume(1h) column (E) values (already text-safe, contains % and spaces)
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("E6").Value = '  -2.74%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("E9").Value = '  -1.88%  '
$ws.Range("E10").Value = '  -4.72%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("E15").Value = '  -2.33%  '
$ws.Range("E16").Value = '  -2.74%  '
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("E20").Value = '  -3.34%  '
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("E22").Value = '  -2.52%  '
$ws.Range("E23").Value = '  -4.89%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("E26").Value = '  -3.88%  '
$ws.Range("E27").Value = '  -1.87%  '
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("E29").Value = '  +0.51%  '
$ws.Range("E30").Value = '  -4.69%  '
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("E33").Value = '  -4.74%  '
$ws.Range("E34").Value = '  -5.94%  '
$ws.Range("E35").Value = '  -4.36%  '
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E38").Value = '  -8.26%  '
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("E40").Value = '  -1.61%  '
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("E43").Value = '  -7.68%  '
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("E47").Value = '  -3.86%  '
$ws.Range("E48").Value = '  +1.26%  '
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("E51").Value = '  -5.64%  '
